$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sighting record is being inserted as row 59, pushing the rows
# currently at 59-62 down to 60-63.
#
# Rows.Insert() fabricates brand-new blank-cell styles for the inserted
# band (and for anything pasted onto it) instead of reusing the existing
# style catalogue, so instead the row band is shifted "by hand":
#   1. Row 62 (last one affected) is cloned onto the brand-new row 63
#      first, while row 62 still holds its original data, using a
#      format-paste followed by a values-paste.
#   2. Rows 61, 60 and 59 are then copied down into 62, 61 and 60 by
#      plain value assignment - each destination row already carries the
#      formatting the shifted content needs (it's identical to what used
#      to be one row above it), so a value-only write leaves every style
#      untouched and correct.
#   3. Row 59 is finally overwritten with the new record.

$ws.Range("A62:I62").Copy()
$ws.Range("A63:I63").PasteSpecial(-4122)
$ws.Range("A62:I62").Copy()
$ws.Range("A63:I63").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$cols = @("A","B","C","D","E","F","G","H","I")
for ($r = 61; $r -ge 59; $r--) {
    $dest = $r + 1
    foreach ($col in $cols) {
        $src = $ws.Range("$col$r")
        $ws.Range("$col$dest").Value = $src.Value2
    }
}

# Write the new record into row 59.
$ws.Range("A59").Value = 45058
$ws.Range("B59").Value = "BUAM"
$ws.Range("C59").Value = 66
$ws.Range("D59").Value = "Jonquière"
$ws.Range("E59").Value = "Saguenay - Lac-Saint-Jean"
$ws.Range("F59").Value = "C"
$ws.Range("G59").Value = "Cote 3"
$ws.Range("H59").Value = "Avec RASY et PSCR"
$ws.Range("I59").Value = "Martin Bertrand"

# F59 (zone climatique = "C") needs the fill used elsewhere for that
# value (e.g. F55) rather than the fill the row it displaced ("A") used.
$ws.Range("F55").Copy()
$ws.Range("F59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D67").Select()
